$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 54.46713274517316
$ws.Range("M4").Value = 39.790391641953093
$ws.Range("M5").Value = 77.295274042005587
$ws.Range("M6").Value = 37.374947253641352
$ws.Range("M7").Value = 31.296339673006749
$ws.Range("M8").Value = 30.431195639155501
$ws.Range("M9").Value = 28.93465931102369
$ws.Range("M10").Value = 171.9960647176099
$ws.Range("M11").Value = 87.903185770664876
$ws.Range("M12").Value = 59.890885753675519
$ws.Range("M13").Value = 63.486656447616639
$ws.Range("M14").Value = 61.617732914796079
$ws.Range("M15").Value = 136.22046273298781
$ws.Range("M16").Value = 462.14463659982363
$ws.Range("M17").Value = 407.9055187250251
$ws.Range("M18").Value = 37.290223656749582
$ws.Range("M19").Value = 37.380204494397667

$ws.Range("T21").Select()
